$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Página1")

# --- Sheet "Página1" updates (bl10.1 content update) ---

# B3: "Natal em família" -> "Passeios, viagens, e Natal em família."
$ws1.Range("B3").Value = "Passeios, viagens, e Natal em família." + [char]10

# F3: append detail about marrying childhood sweetheart
$ws1.Range("F3").Value = "Contratado para trabalhar no RH de uma multinacional, morei dois anos em outra cidade, me formei em administração, nascimento do meu filho e casei com a meu amor platônico da época de escola. "

# G3: fix typo "Decidir" -> "Decidi"
$ws1.Range("G3").Value = "Decidi fazer transição de carreira e começar a estudar na Trybe."

# E4: new negative moment about the car
$ws1.Range("E4").Value = "Comprei um carro parcelado que bateu o motor e fiquei mais com ele parado, e pagando as parcelas, do que com ele andando. "

# Update row 3 height (wrapped text grew)
$ws1.Rows.Item(3).RowHeight = 124.6

# Update active cell selection on sheet 1
$ws1.Range("G5").Select()

$wb.Save()
